# Generate Report for handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) values on the zh-cn and de-de
# sheets to reflect a newly re-run report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-16 10:38:16"
$wsZhCn.Range("D3").Value = "2016-02-16 10:38:16"
$wsZhCn.Range("G2").Value = "2016-02-16 10:39:12"
$wsZhCn.Range("G3").Value = "2016-02-16 10:39:12"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-16 10:38:30"
$wsDeDe.Range("D3").Value = "2016-02-16 10:38:30"
$wsDeDe.Range("G2").Value = "2016-02-16 10:39:39"
$wsDeDe.Range("G3").Value = "2016-02-16 10:39:39"
